# Recalculated profit-margin figures (currentAveragePrice / LevePrice / LeveProfit
# columns H-N of each Table_<Role> sheet) refreshed by the scheduled market-data
# runner. Writes the new values cell-by-cell per sheet; two rows (ARM!N32 and
# CUL!M131) previously had no value in that column, so assigning .Value adds the cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 175.36363  # H18
$ws.Cells.Item(18, 9).Value = 150.38095  # I18
$ws.Cells.Item(18, 11).Value = 150.38095  # K18
$ws.Cells.Item(18, 13).Value = 133.61905  # M18
$ws.Cells.Item(40, 8).Value = 1618.125  # H40
$ws.Cells.Item(40, 9).Value = 1038.75  # I40
$ws.Cells.Item(40, 10).Value = 2197.5  # J40
$ws.Cells.Item(40, 11).Value = 1038.75  # K40
$ws.Cells.Item(40, 12).Value = 2197.5  # L40
$ws.Cells.Item(40, 13).Value = -863.75  # M40
$ws.Cells.Item(40, 14).Value = -2547.5  # N40
$ws.Cells.Item(112, 8).Value = 1025.8485  # H112
$ws.Cells.Item(112, 10).Value = 1056.5483  # J112
$ws.Cells.Item(112, 12).Value = 3169.6449  # L112
$ws.Cells.Item(112, 14).Value = -5385.644899999999  # N112
$ws.Cells.Item(116, 8).Value = 4291.4546  # H116
$ws.Cells.Item(116, 9).Value = 1898  # I116
$ws.Cells.Item(116, 11).Value = 1898  # K116
$ws.Cells.Item(116, 13).Value = 1544  # M116
$ws.Cells.Item(129, 8).Value = 182641.44  # H129
$ws.Cells.Item(129, 10).Value = 193157.28  # J129
$ws.Cells.Item(129, 12).Value = 579471.84  # L129
$ws.Cells.Item(129, 14).Value = -589471.84  # N129
$ws.Cells.Item(137, 8).Value = 74971.42999999999  # H137
$ws.Cells.Item(137, 9).Value = 4712.375  # I137
$ws.Cells.Item(137, 11).Value = 14137.125  # K137
$ws.Cells.Item(137, 13).Value = -11587.125  # M137
$ws.Cells.Item(138, 8).Value = 1802.096  # H138
$ws.Cells.Item(138, 9).Value = 1138.3889  # I138
$ws.Cells.Item(138, 10).Value = 2447.8647  # J138
$ws.Cells.Item(138, 11).Value = 3415.1667  # K138
$ws.Cells.Item(138, 12).Value = 7343.5941  # L138
$ws.Cells.Item(138, 13).Value = 1724.8333  # M138
$ws.Cells.Item(138, 14).Value = -17623.5941  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 21960.76  # H32
$ws.Cells.Item(32, 9).Value = 22394.652  # I32
$ws.Cells.Item(32, 10).Value = 700  # J32
$ws.Cells.Item(32, 11).Value = 22394.652  # K32
$ws.Cells.Item(32, 12).Value = 700  # L32
$ws.Cells.Item(32, 13).Value = -22107.652  # M32
$ws.Cells.Item(32, 14).Value = -1274  # N32
$ws.Cells.Item(61, 8).Value = 2261.0605  # H61
$ws.Cells.Item(61, 9).Value = 1825.4  # I61
$ws.Cells.Item(61, 10).Value = 3622.5  # J61
$ws.Cells.Item(61, 11).Value = 1825.4  # K61
$ws.Cells.Item(61, 12).Value = 3622.5  # L61
$ws.Cells.Item(61, 13).Value = -1613.4  # M61
$ws.Cells.Item(61, 14).Value = -4046.5  # N61
$ws.Cells.Item(74, 8).Value = 43480690  # H74
$ws.Cells.Item(74, 9).Value = 50002630  # I74
$ws.Cells.Item(74, 10).Value = 1066.6666  # J74
$ws.Cells.Item(74, 11).Value = 50002630  # K74
$ws.Cells.Item(74, 12).Value = 1066.6666  # L74
$ws.Cells.Item(74, 13).Value = -50001756  # M74
$ws.Cells.Item(74, 14).Value = -2814.6666  # N74
$ws.Cells.Item(77, 8).Value = 43480690  # H77
$ws.Cells.Item(77, 9).Value = 50002630  # I77
$ws.Cells.Item(77, 10).Value = 1066.6666  # J77
$ws.Cells.Item(77, 11).Value = 250013150  # K77
$ws.Cells.Item(77, 12).Value = 5333.333000000001  # L77
$ws.Cells.Item(77, 13).Value = -250008782  # M77
$ws.Cells.Item(77, 14).Value = -14069.333  # N77
$ws.Cells.Item(97, 8).Value = 1919.875  # H97
$ws.Cells.Item(97, 9).Value = 2151.2856  # I97
$ws.Cells.Item(97, 10).Value = 300  # J97
$ws.Cells.Item(97, 11).Value = 2151.2856  # K97
$ws.Cells.Item(97, 12).Value = 300  # L97
$ws.Cells.Item(97, 13).Value = -1655.2856  # M97
$ws.Cells.Item(97, 14).Value = -1292  # N97
$ws.Cells.Item(102, 8).Value = 1683.3  # H102
$ws.Cells.Item(102, 9).Value = 1472.1666  # I102
$ws.Cells.Item(102, 10).Value = 2000  # J102
$ws.Cells.Item(102, 11).Value = 1472.1666  # K102
$ws.Cells.Item(102, 12).Value = 2000  # L102
$ws.Cells.Item(102, 13).Value = 149.8334  # M102
$ws.Cells.Item(102, 14).Value = -5244  # N102
$ws.Cells.Item(132, 8).Value = 10822.926  # H132
$ws.Cells.Item(132, 9).Value = 1383.2195  # I132
$ws.Cells.Item(132, 11).Value = 4149.6585  # K132
$ws.Cells.Item(132, 13).Value = -1619.6585  # M132
$ws.Cells.Item(135, 8).Value = 26475  # H135
$ws.Cells.Item(135, 10).Value = 26475  # J135
$ws.Cells.Item(135, 12).Value = 26475  # L135
$ws.Cells.Item(135, 14).Value = -36615  # N135
$ws.Cells.Item(136, 8).Value = 2261.0605  # H136
$ws.Cells.Item(136, 9).Value = 1825.4  # I136
$ws.Cells.Item(136, 10).Value = 3622.5  # J136
$ws.Cells.Item(136, 11).Value = 5476.200000000001  # K136
$ws.Cells.Item(136, 12).Value = 10867.5  # L136
$ws.Cells.Item(136, 13).Value = -2926.200000000001  # M136
$ws.Cells.Item(136, 14).Value = -15967.5  # N136
$ws.Cells.Item(139, 8).Value = 39997.69  # H139
$ws.Cells.Item(139, 10).Value = 39997.69  # J139
$ws.Cells.Item(139, 12).Value = 39997.69  # L139
$ws.Cells.Item(139, 14).Value = -50277.69  # N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 27938.75  # H81
$ws.Cells.Item(81, 10).Value = 27938.75  # J81
$ws.Cells.Item(81, 12).Value = 27938.75  # L81
$ws.Cells.Item(81, 14).Value = -30060.75  # N81
$ws.Cells.Item(84, 8).Value = 27938.75  # H84
$ws.Cells.Item(84, 10).Value = 27938.75  # J84
$ws.Cells.Item(84, 12).Value = 83816.25  # L84
$ws.Cells.Item(84, 14).Value = -94424.25  # N84
$ws.Cells.Item(134, 8).Value = 26888.785  # H134
$ws.Cells.Item(134, 9).Value = 31966.543  # I134
$ws.Cells.Item(134, 11).Value = 95899.629  # K134
$ws.Cells.Item(134, 13).Value = -93364.629  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 12229.607  # H31
$ws.Cells.Item(31, 9).Value = 14157.782  # I31
$ws.Cells.Item(31, 11).Value = 14157.782  # K31
$ws.Cells.Item(31, 13).Value = -13862.782  # M31
$ws.Cells.Item(34, 8).Value = 12229.607  # H34
$ws.Cells.Item(34, 9).Value = 14157.782  # I34
$ws.Cells.Item(34, 11).Value = 14157.782  # K34
$ws.Cells.Item(34, 13).Value = -13955.782  # M34
$ws.Cells.Item(58, 8).Value = 15753.059  # H58
$ws.Cells.Item(58, 9).Value = 1113.6086  # I58
$ws.Cells.Item(58, 10).Value = 46362.816  # J58
$ws.Cells.Item(58, 11).Value = 1113.6086  # K58
$ws.Cells.Item(58, 12).Value = 46362.816  # L58
$ws.Cells.Item(58, 13).Value = -910.6086  # M58
$ws.Cells.Item(58, 14).Value = -46768.816  # N58
$ws.Cells.Item(60, 8).Value = 10173.685  # H60
$ws.Cells.Item(60, 9).Value = 2999.3333  # I60
$ws.Cells.Item(60, 10).Value = 11518.875  # J60
$ws.Cells.Item(60, 11).Value = 2999.3333  # K60
$ws.Cells.Item(60, 12).Value = 11518.875  # L60
$ws.Cells.Item(60, 13).Value = -2488.3333  # M60
$ws.Cells.Item(60, 14).Value = -12540.875  # N60
$ws.Cells.Item(62, 8).Value = 166670600  # H62
$ws.Cells.Item(62, 10).Value = 5168.3335  # J62
$ws.Cells.Item(62, 12).Value = 5168.3335  # L62
$ws.Cells.Item(62, 14).Value = -6416.3335  # N62
$ws.Cells.Item(65, 8).Value = 166670600  # H65
$ws.Cells.Item(65, 10).Value = 5168.3335  # J65
$ws.Cells.Item(65, 12).Value = 25841.6675  # L65
$ws.Cells.Item(65, 14).Value = -32081.6675  # N65
$ws.Cells.Item(99, 8).Value = 14289570  # H99
$ws.Cells.Item(99, 9).Value = 3118  # I99
$ws.Cells.Item(99, 11).Value = 3118  # K99
$ws.Cells.Item(99, 13).Value = -1620  # M99
$ws.Cells.Item(126, 8).Value = 14289570  # H126
$ws.Cells.Item(126, 9).Value = 3118  # I126
$ws.Cells.Item(126, 11).Value = 9354  # K126
$ws.Cells.Item(126, 13).Value = -6884  # M126
$ws.Cells.Item(132, 8).Value = 13557.907  # H132
$ws.Cells.Item(132, 9).Value = 17270.871  # I132
$ws.Cells.Item(132, 11).Value = 51812.613  # K132
$ws.Cells.Item(132, 13).Value = -49282.613  # M132
$ws.Cells.Item(134, 8).Value = 1081.4565  # H134
$ws.Cells.Item(134, 9).Value = 951.8333  # I134
$ws.Cells.Item(134, 10).Value = 1164.7858  # J134
$ws.Cells.Item(134, 11).Value = 2855.4999  # K134
$ws.Cells.Item(134, 12).Value = 3494.3574  # L134
$ws.Cells.Item(134, 13).Value = -320.4998999999998  # M134
$ws.Cells.Item(134, 14).Value = -8564.357400000001  # N134
$ws.Cells.Item(136, 8).Value = 15753.059  # H136
$ws.Cells.Item(136, 9).Value = 1113.6086  # I136
$ws.Cells.Item(136, 10).Value = 46362.816  # J136
$ws.Cells.Item(136, 11).Value = 3340.8258  # K136
$ws.Cells.Item(136, 12).Value = 139088.448  # L136
$ws.Cells.Item(136, 13).Value = -790.8258000000001  # M136
$ws.Cells.Item(136, 14).Value = -144188.448  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(60, 8).Value = 1050  # H60
$ws.Cells.Item(113, 8).Value = 2602.3044  # H113
$ws.Cells.Item(113, 9).Value = 3143.0557  # I113
$ws.Cells.Item(113, 11).Value = 9429.167099999999  # K113
$ws.Cells.Item(113, 13).Value = -7259.167099999999  # M113
$ws.Cells.Item(122, 8).Value = 852.7  # H122
$ws.Cells.Item(122, 9).Value = 412.25  # I122
$ws.Cells.Item(122, 10).Value = 1146.3334  # J122
$ws.Cells.Item(122, 11).Value = 3710.25  # K122
$ws.Cells.Item(122, 12).Value = 10317.0006  # L122
$ws.Cells.Item(122, 13).Value = -1260.25  # M122
$ws.Cells.Item(122, 14).Value = -15217.0006  # N122
$ws.Cells.Item(131, 8).Value = 739.7  # H131
$ws.Cells.Item(131, 9).Value = 300  # I131
$ws.Cells.Item(131, 10).Value = 744.1414  # J131
$ws.Cells.Item(131, 11).Value = 900  # K131
$ws.Cells.Item(131, 12).Value = 2232.4242  # L131
$ws.Cells.Item(131, 13).Value = 4140  # M131
$ws.Cells.Item(131, 14).Value = -12312.4242  # N131
$ws.Cells.Item(132, 8).Value = 1244.6  # H132
$ws.Cells.Item(132, 9).Value = 1059.8572  # I132
$ws.Cells.Item(132, 10).Value = 1406.25  # J132
$ws.Cells.Item(132, 11).Value = 9538.7148  # K132
$ws.Cells.Item(132, 12).Value = 12656.25  # L132
$ws.Cells.Item(132, 13).Value = -7008.7148  # M132
$ws.Cells.Item(132, 14).Value = -17716.25  # N132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 417.2  # H97
$ws.Cells.Item(97, 9).Value = 375.57144  # I97
$ws.Cells.Item(97, 11).Value = 375.57144  # K97
$ws.Cells.Item(97, 13).Value = 120.42856  # M97
$ws.Cells.Item(102, 8).Value = 25002322  # H102
$ws.Cells.Item(102, 9).Value = 29414324  # I102
$ws.Cells.Item(102, 10).Value = 978  # J102
$ws.Cells.Item(102, 11).Value = 29414324  # K102
$ws.Cells.Item(102, 12).Value = 978  # L102
$ws.Cells.Item(102, 13).Value = -29412702  # M102
$ws.Cells.Item(102, 14).Value = -4222  # N102
$ws.Cells.Item(126, 8).Value = 3410.1135  # H126
$ws.Cells.Item(126, 9).Value = 2601.5  # I126
$ws.Cells.Item(126, 11).Value = 7804.5  # K126
$ws.Cells.Item(126, 13).Value = -5334.5  # M126

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1114  # H46
$ws.Cells.Item(46, 9).Value = 830.3333  # I46
$ws.Cells.Item(46, 10).Value = 1235.5714  # J46
$ws.Cells.Item(46, 11).Value = 830.3333  # K46
$ws.Cells.Item(46, 12).Value = 1235.5714  # L46
$ws.Cells.Item(46, 13).Value = -642.3333  # M46
$ws.Cells.Item(46, 14).Value = -1611.5714  # N46
$ws.Cells.Item(82, 8).Value = 2397.9  # H82
$ws.Cells.Item(82, 9).Value = 2476  # I82
$ws.Cells.Item(82, 10).Value = 2319.8  # J82
$ws.Cells.Item(82, 11).Value = 2476  # K82
$ws.Cells.Item(82, 12).Value = 2319.8  # L82
$ws.Cells.Item(82, 13).Value = -2115  # M82
$ws.Cells.Item(82, 14).Value = -3041.8  # N82
$ws.Cells.Item(85, 8).Value = 2397.9  # H85
$ws.Cells.Item(85, 9).Value = 2476  # I85
$ws.Cells.Item(85, 10).Value = 2319.8  # J85
$ws.Cells.Item(85, 11).Value = 2476  # K85
$ws.Cells.Item(85, 12).Value = 2319.8  # L85
$ws.Cells.Item(85, 13).Value = -1228  # M85
$ws.Cells.Item(85, 14).Value = -4815.8  # N85
$ws.Cells.Item(100, 8).Value = 2306.5715  # H100
$ws.Cells.Item(100, 9).Value = 1998.75  # I100
$ws.Cells.Item(100, 11).Value = 1998.75  # K100
$ws.Cells.Item(100, 13).Value = -1457.75  # M100
$ws.Cells.Item(122, 8).Value = 1035097.2  # H122
$ws.Cells.Item(122, 9).Value = 1510842.5  # I122
$ws.Cells.Item(122, 10).Value = 4315.8335  # J122
$ws.Cells.Item(122, 11).Value = 4532527.5  # K122
$ws.Cells.Item(122, 12).Value = 12947.5005  # L122
$ws.Cells.Item(122, 13).Value = -4530077.5  # M122
$ws.Cells.Item(122, 14).Value = -17847.5005  # N122
$ws.Cells.Item(132, 8).Value = 2315.25  # H132
$ws.Cells.Item(132, 9).Value = 1571.2  # I132
$ws.Cells.Item(132, 10).Value = 3555.3333  # J132
$ws.Cells.Item(132, 11).Value = 4713.6  # K132
$ws.Cells.Item(132, 12).Value = 10665.9999  # L132
$ws.Cells.Item(132, 13).Value = -2183.6  # M132
$ws.Cells.Item(132, 14).Value = -15725.9999  # N132
$ws.Cells.Item(136, 8).Value = 24354.637  # H136
$ws.Cells.Item(136, 9).Value = 36807.285  # I136
$ws.Cells.Item(136, 10).Value = 2562.5  # J136
$ws.Cells.Item(136, 11).Value = 110421.855  # K136
$ws.Cells.Item(136, 12).Value = 7687.5  # L136
$ws.Cells.Item(136, 13).Value = -107871.855  # M136
$ws.Cells.Item(136, 14).Value = -12787.5  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1283.9231  # H126
$ws.Cells.Item(126, 9).Value = 1431.9333  # I126
$ws.Cells.Item(126, 11).Value = 4295.7999  # K126
$ws.Cells.Item(126, 13).Value = -1825.7999  # M126
$ws.Cells.Item(132, 8).Value = 1594.2106  # H132
$ws.Cells.Item(132, 9).Value = 965.5833  # I132
$ws.Cells.Item(132, 10).Value = 2671.8572  # J132
$ws.Cells.Item(132, 11).Value = 2896.7499  # K132
$ws.Cells.Item(132, 12).Value = 8015.571599999999  # L132
$ws.Cells.Item(132, 13).Value = -366.7498999999998  # M132
$ws.Cells.Item(132, 14).Value = -13075.5716  # N132
$ws.Cells.Item(136, 8).Value = 19609404  # H136
$ws.Cells.Item(136, 9).Value = 28572812  # I136
$ws.Cells.Item(136, 11).Value = 85718436  # K136
$ws.Cells.Item(136, 13).Value = -85715886  # M136
